$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C (Förändrad / "changed date") for all existing data rows (2-359)
#    from 45190 to 45192 (2023-09-23), keeping the existing date style.
$rng = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item(359, 3))
$rng.Value = 45192

# 2. Row 359 picks up an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(359).RowHeight = 15

# 3. Append the new record as row 360.
$ws.Cells.Item(360, 1).Value = "A 44629-2023"

$ws.Cells.Item(360, 2).Value = 45189
$ws.Cells.Item(360, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(360, 3).Value = 45192
$ws.Cells.Item(360, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(360, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(360, 5).Value = "BJURHOLM"

$ws.Cells.Item(360, 7).Value = 15.6
$ws.Cells.Item(360, 8).Value = 0
$ws.Cells.Item(360, 9).Value = 0
$ws.Cells.Item(360, 10).Value = 0
$ws.Cells.Item(360, 11).Value = 0
$ws.Cells.Item(360, 12).Value = 0
$ws.Cells.Item(360, 13).Value = 0
$ws.Cells.Item(360, 14).Value = 0
$ws.Cells.Item(360, 15).Value = 0
$ws.Cells.Item(360, 16).Value = 0
$ws.Cells.Item(360, 17).Value = 0

$ws.Cells.Item(360, 18).WrapText = $true
